# Generate Report for Handoff
#
# The localization-status report moves from "In Translation" to
# "Ready for handoff" and its timestamps are refreshed to the moment the
# handoff package was produced. Widening the "Status" columns mirrors
# Excel's column auto-fit kicking in for the longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ----------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps --------------------------------------------------
# Overview's "Latest HO Xliff Generate Date" and de-de's "Latest Handoff
# Datetime" share one timestamp; zh-cn's handoff datetime gets its own.
$wsOverview.Range("G2").Value = "2016-08-16 00:35:55"
$wsDeDe.Range("H2").Value     = "2016-08-16 00:35:55"
$wsZhCn.Range("H2").Value     = "2016-08-16 00:35:50"

# --- Column widths: Status columns widen to fit "Ready for handoff" -------
$newWidth = 98 / 6
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth     = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth     = $newWidth
